$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Poongodi.R"
$ws.Range("C13").Value = "02.01.2018"
$ws.Range("D13").Value = "8.30 to 4.30"
$ws.Range("E13").Value = "project flow,API documentation"
$ws.Range("F13").Value = "completed"

$ws.Range("C19").Value = ","
